$d = $word.ActiveDocument

# Locate the target phrase "1- Given the following sets: " so we don't rely
# on hard-coded character offsets.
$findRange = $d.Content
$found = $findRange.Find.Execute("1- Given the following sets:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $findRange.Start   # position of "1"
$dashPos = $start + 1       # position of "-"

# Step 1: split "1" away from "- Given the following sets: " by toggling a
# throw-away formatting change on "1" (forces the host to keep it as its own
# run) and then replace the "-" character with ".".
$r1 = $d.Range($start, $start + 1)        # "1"
$rDash = $d.Range($dashPos, $dashPos + 1) # "-"

$r1.Bold = 1
$rDash.Text = "."
$r1.Bold = 0

# Step 2: split the new "." away from the remaining " Given the following
# sets: " text using the same trick, now that the "." sits in its own run
# boundary at $dashPos .. $dashPos+1.
$rDot = $d.Range($dashPos, $dashPos + 1)  # "."
$rDot.Bold = 1
$rDot.Bold = 0
